$wb = $excel.ActiveWorkbook

# Sheet 1: D.chrysitis — fix trailing-space "Wiejce " -> "Wiejce"
$ws1 = $wb.Worksheets.Item("D.chrysitis")
$ws1.Range("F18").Value = "Wiejce"
$ws1.Range("F22").Value = "Wiejce"
$ws1.Range("F42").Value = "Wiejce"

# Sheet 2: D.stenochrysis — fix trailing-space localities
$ws2 = $wb.Worksheets.Item("D.stenochrysis")

$ws2.Range("F9").Value = "Czelin"
$ws2.Range("F11").Value = "Czelin"
$ws2.Range("F13").Value = "Czelin"

$ws2.Range("F32").Value = "Świnoujście"
$ws2.Range("F33").Value = "Świnoujście"
$ws2.Range("F34").Value = "Świnoujście"
$ws2.Range("F35").Value = "Świnoujście"
$ws2.Range("F36").Value = "Świnoujście"
$ws2.Range("F37").Value = "Świnoujście"

$ws2.Range("F38").Value = "Woźniki"
$ws2.Range("F39").Value = "Woźniki"
$ws2.Range("F40").Value = "Woźniki"
$ws2.Range("F41").Value = "Woźniki"
$ws2.Range("F42").Value = "Woźniki"
$ws2.Range("F43").Value = "Woźniki"
